# Remove the stray "Sheet" audit row from the optimization_parameters sheet.
# (This leftover row — A16="Sheet", B16=3, C16=4 — was a test artifact; the
# author's commit "wrapping up test file audit" deletes the whole row 16,
# which shifts the former row 17 ("simulation_timepoints" ... ) up into its
# place.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

$ws.Activate()
$ws.Rows("16:16").Select()
$ws.Rows("16:16").Delete()
